$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (matches workbook.xml sheet name change)
$ws.Name = "hate-crimes-noncampus-virginia-"

# Remove the old title row (row 1 held "Hate Crimes - Noncampus" spanning the sheet);
# the former header row (row 2) becomes row 1 and every data row shifts up by one.
$ws.Rows.Item(1).Delete()

# Normalize header capitalization/wording in the new row 1
$ws.Range("A1").Value = "Survey Year"
$ws.Range("B1").Value = "UnitID"
$ws.Range("C1").Value = "Institution Name"
$ws.Range("G1").Value = "Murder/Non-Negligent Manslaughter"
$ws.Range("H1").Value = "Murder/Non-Negligent Manslaughter - Race"
$ws.Range("I1").Value = "Murder/Non-Negligent Manslaughter - Religion"
$ws.Range("J1").Value = "Murder/Non-Negligent Manslaughter - Sexual Orientation"
$ws.Range("K1").Value = "Murder/Non-Negligent Manslaughter - Gender"
$ws.Range("L1").Value = "Murder/Non-Negligent Manslaughter - Disability"
$ws.Range("M1").Value = "Murder/Non-Negligent Manslaughter - Ethnicity/National Origin"
$ws.Range("N1").Value = "Negligent Manslaughter"
$ws.Range("O1").Value = "Negligent Manslaughter - Race"
$ws.Range("P1").Value = "Negligent Manslaughter - Religion"
$ws.Range("Q1").Value = "Negligent Manslaughter - Sexual Orientation"
$ws.Range("R1").Value = "Negligent Manslaughter - Gender"
$ws.Range("S1").Value = "Negligent Manslaughter - Disability"
$ws.Range("T1").Value = "Negligent Manslaughter - Ethnicity/National Origin"
$ws.Range("U1").Value = "Sex Offenses - Forcible"
$ws.Range("V1").Value = "Sex Offenses - Forcible - Race"
$ws.Range("W1").Value = "Sex Offenses - Forcible - Religion"
$ws.Range("X1").Value = "Sex Offenses - Forcible - Sexual Orientation"
$ws.Range("Y1").Value = "Sex Offenses - Forcible - Gender"
$ws.Range("Z1").Value = "Sex Offenses - Forcible - Disability"
$ws.Range("AA1").Value = "Sex Offenses - Forcible - Ethnicity/National Origin"
$ws.Range("AB1").Value = "Sex Offenses - Non-Forcible"
$ws.Range("AC1").Value = "Sex Offenses - Non-Forcible -Race"
$ws.Range("AD1").Value = "Sex Offenses - Non-Forcible - Religion"
$ws.Range("AE1").Value = "Sex Offenses - Non-Forcible - Sexual Orientation"
$ws.Range("AF1").Value = "Sex Offenses - Non-Forcible - Gender"
$ws.Range("AG1").Value = "Sex Offenses - Non-Forcible - Disability"
$ws.Range("AH1").Value = "Sex Offenses - Non-Forcible - Ethnicity/National Origin"
$ws.Range("AL1").Value = "Robbery - Sexual Orientation"
$ws.Range("AO1").Value = "Robbery - Ethnicity/National Origin"
$ws.Range("AP1").Value = "Aggravated Assault"
$ws.Range("AQ1").Value = "Aggravated Assault - Race"
$ws.Range("AR1").Value = "Aggravated Assault - Religion"
$ws.Range("AS1").Value = "Aggravated Assault - Sexual Orientation"
$ws.Range("AT1").Value = "Aggravated Assault - Gender"
$ws.Range("AU1").Value = "Aggravated Assault - Disability"
$ws.Range("AV1").Value = "Aggravated Assault - Ethnicity/National Origin"
$ws.Range("AZ1").Value = "Burglary - Sexual Orientation"
$ws.Range("BC1").Value = "Burglary - Ethnicity/National Origin"
$ws.Range("BD1").Value = "Motor Vehicle Theft"
$ws.Range("BE1").Value = "Motor Vehicle Theft - Race"
$ws.Range("BF1").Value = "Motor Vehicle Theft - Religion"
$ws.Range("BG1").Value = "Motor Vehicle Theft - Sexual Orientation"
$ws.Range("BH1").Value = "Motor Vehicle Theft - Gender"
$ws.Range("BI1").Value = "Motor Vehicle Theft - Disability"
$ws.Range("BJ1").Value = "Motor Vehicle Theft - Ethnicity/National Origin"
$ws.Range("BN1").Value = "Arson - Sexual Orientation"
$ws.Range("BQ1").Value = "Arson - Ethnicity/National Origin"
$ws.Range("BR1").Value = "Simple Assault"
$ws.Range("BS1").Value = "Simple Assault - Race"
$ws.Range("BT1").Value = "Simple Assault - Religion"
$ws.Range("BU1").Value = "Simple Assault - Sexual Orientation"
$ws.Range("BV1").Value = "Simple Assault - Gender"
$ws.Range("BW1").Value = "Simple Assault - Disability"
$ws.Range("BX1").Value = "Simple Assault - Ethnicity/National Origin"
$ws.Range("BY1").Value = "Larceny-Theft"
$ws.Range("BZ1").Value = "Larceny-Theft - Race"
$ws.Range("CA1").Value = "Larceny-Theft - Religion"
$ws.Range("CB1").Value = "Larceny-Theft - Sexual Orientation"
$ws.Range("CC1").Value = "Larceny-Theft - Gender"
$ws.Range("CD1").Value = "Larceny-Theft - Disability"
$ws.Range("CE1").Value = "Larceny-Theft - Ethnicity/National Origin"
$ws.Range("CI1").Value = "Intimidation - Sexual Orientation"
$ws.Range("CL1").Value = "Intimidation - Ethnicity/National Origin"
$ws.Range("CM1").Value = "Destruction/Damage/Vandalism of Property"
$ws.Range("CN1").Value = "Destruction/Damage/Vandalism of Property - Race"
$ws.Range("CO1").Value = "Destruction/Damage/Vandalism of Property - Religion"
$ws.Range("CP1").Value = "Destruction/Damage/Vandalism of Property - Sexual Orientation"
$ws.Range("CQ1").Value = "Destruction/Damage/Vandalism of Property - Gender"
$ws.Range("CR1").Value = "Destruction/Damage/Vandalism of Property - Disability"
$ws.Range("CS1").Value = "Destruction/Damage/Vandalism of Property - Ethnicity/National Origin"
